# Weekly update: insert 5 new price rows for the week of 2023-12-05
# (serial date 45265) into the "Fruta, Terminal La Palmera de La Serena -
# Cereza" sheet, pushing the existing data rows down by five positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert five blank rows above the current row 420; this shifts the old
# rows 420:496 down to 425:501 and grows the sheet from T496 to T501.
$ws.Rows("420:424").Insert()

# Shared attribute values for every new row in this weekly block.
$mercadoId = 8
$mercado = "Terminal La Palmera de La Serena"
$region = "Coquimbo"
$fecha = 45265
$codreg = 4
$tipo = "Fruta"
$productoId = 100103
$producto = "Frutos de hueso (carozo)"
$categoriaId = 100103001
$categoria = "Cereza"
$unidad = "`$/bandeja 10 kilos"
$kgUnidad = 10

# Row 420: Lapins / Primera
$ws.Cells.Item(420, 1).Value = $mercadoId
$ws.Cells.Item(420, 2).Value = $mercado
$ws.Cells.Item(420, 3).Value = $region
$ws.Cells.Item(420, 4).Value = $fecha
$ws.Cells.Item(420, 5).Value = $codreg
$ws.Cells.Item(420, 6).Value = $tipo
$ws.Cells.Item(420, 7).Value = $productoId
$ws.Cells.Item(420, 8).Value = $producto
$ws.Cells.Item(420, 9).Value = $categoriaId
$ws.Cells.Item(420, 10).Value = $categoria
$ws.Cells.Item(420, 11).Value = "Lapins"
$ws.Cells.Item(420, 12).Value = "Primera"
$ws.Cells.Item(420, 13).Value = 400
$ws.Cells.Item(420, 14).Value = 14000
$ws.Cells.Item(420, 15).Value = 15000
$ws.Cells.Item(420, 16).Value = 14500
$ws.Cells.Item(420, 17).Value = $unidad
$ws.Cells.Item(420, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(420, 19).Value = 1450
$ws.Cells.Item(420, 20).Value = $kgUnidad

# Row 421: Lapins / Segunda
$ws.Cells.Item(421, 1).Value = $mercadoId
$ws.Cells.Item(421, 2).Value = $mercado
$ws.Cells.Item(421, 3).Value = $region
$ws.Cells.Item(421, 4).Value = $fecha
$ws.Cells.Item(421, 5).Value = $codreg
$ws.Cells.Item(421, 6).Value = $tipo
$ws.Cells.Item(421, 7).Value = $productoId
$ws.Cells.Item(421, 8).Value = $producto
$ws.Cells.Item(421, 9).Value = $categoriaId
$ws.Cells.Item(421, 10).Value = $categoria
$ws.Cells.Item(421, 11).Value = "Lapins"
$ws.Cells.Item(421, 12).Value = "Segunda"
$ws.Cells.Item(421, 13).Value = 300
$ws.Cells.Item(421, 14).Value = 12000
$ws.Cells.Item(421, 15).Value = 13000
$ws.Cells.Item(421, 16).Value = 12500
$ws.Cells.Item(421, 17).Value = $unidad
$ws.Cells.Item(421, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(421, 19).Value = 1250
$ws.Cells.Item(421, 20).Value = $kgUnidad

# Row 422: Rainier / Primera
$ws.Cells.Item(422, 1).Value = $mercadoId
$ws.Cells.Item(422, 2).Value = $mercado
$ws.Cells.Item(422, 3).Value = $region
$ws.Cells.Item(422, 4).Value = $fecha
$ws.Cells.Item(422, 5).Value = $codreg
$ws.Cells.Item(422, 6).Value = $tipo
$ws.Cells.Item(422, 7).Value = $productoId
$ws.Cells.Item(422, 8).Value = $producto
$ws.Cells.Item(422, 9).Value = $categoriaId
$ws.Cells.Item(422, 10).Value = $categoria
$ws.Cells.Item(422, 11).Value = "Rainier"
$ws.Cells.Item(422, 12).Value = "Primera"
$ws.Cells.Item(422, 13).Value = 360
$ws.Cells.Item(422, 14).Value = 14000
$ws.Cells.Item(422, 15).Value = 15000
$ws.Cells.Item(422, 16).Value = 14500
$ws.Cells.Item(422, 17).Value = $unidad
$ws.Cells.Item(422, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(422, 19).Value = 1450
$ws.Cells.Item(422, 20).Value = $kgUnidad

# Row 423: Santina / Primera
$ws.Cells.Item(423, 1).Value = $mercadoId
$ws.Cells.Item(423, 2).Value = $mercado
$ws.Cells.Item(423, 3).Value = $region
$ws.Cells.Item(423, 4).Value = $fecha
$ws.Cells.Item(423, 5).Value = $codreg
$ws.Cells.Item(423, 6).Value = $tipo
$ws.Cells.Item(423, 7).Value = $productoId
$ws.Cells.Item(423, 8).Value = $producto
$ws.Cells.Item(423, 9).Value = $categoriaId
$ws.Cells.Item(423, 10).Value = $categoria
$ws.Cells.Item(423, 11).Value = "Santina"
$ws.Cells.Item(423, 12).Value = "Primera"
$ws.Cells.Item(423, 13).Value = 340
$ws.Cells.Item(423, 14).Value = 14500
$ws.Cells.Item(423, 15).Value = 15000
$ws.Cells.Item(423, 16).Value = 14750
$ws.Cells.Item(423, 17).Value = $unidad
$ws.Cells.Item(423, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(423, 19).Value = 1475
$ws.Cells.Item(423, 20).Value = $kgUnidad

# Row 424: Santina / Segunda
$ws.Cells.Item(424, 1).Value = $mercadoId
$ws.Cells.Item(424, 2).Value = $mercado
$ws.Cells.Item(424, 3).Value = $region
$ws.Cells.Item(424, 4).Value = $fecha
$ws.Cells.Item(424, 5).Value = $codreg
$ws.Cells.Item(424, 6).Value = $tipo
$ws.Cells.Item(424, 7).Value = $productoId
$ws.Cells.Item(424, 8).Value = $producto
$ws.Cells.Item(424, 9).Value = $categoriaId
$ws.Cells.Item(424, 10).Value = $categoria
$ws.Cells.Item(424, 11).Value = "Santina"
$ws.Cells.Item(424, 12).Value = "Segunda"
$ws.Cells.Item(424, 13).Value = 300
$ws.Cells.Item(424, 14).Value = 11500
$ws.Cells.Item(424, 15).Value = 12000
$ws.Cells.Item(424, 16).Value = 11750
$ws.Cells.Item(424, 17).Value = $unidad
$ws.Cells.Item(424, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(424, 19).Value = 1175
$ws.Cells.Item(424, 20).Value = $kgUnidad
